$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Numeric-looking Price strings must be written as TEXT (not auto-converted
# to a Double, which would corrupt trailing zeros / precision), so for those
# cells we temporarily force a Text number format, assign the value, then
# restore the cell style so no visual/style diff is introduced.

$ws.Range("D2").Value = '35.354.45'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.904.40'
$ws.Range("E3").Value = '  +2.32%  '

$ws.Range("E4").Value = '  -0.43%  '

$d5 = $ws.Range("D5")
$d5.NumberFormat = "@"
$d5.Value = '246.25'
$d5.Style = "Normal"
$ws.Range("E5").Value = '  +2.94%  '

$ws.Range("E6").Value = '  +6.21%  '

$ws.Range("E7").Value = '  -0.42%  '

$d8 = $ws.Range("D8")
$d8.NumberFormat = "@"
$d8.Value = '41.50'
$d8.Style = "Normal"
$ws.Range("E8").Value = '  -2.12%  '

$ws.Range("E9").Value = '  +5.50%  '

$d10 = $ws.Range("D10")
$d10.NumberFormat = "@"
$d10.Value = '53.03'
$d10.Style = "Normal"
$ws.Range("E10").Value = '  +12.98%  '

$d11 = $ws.Range("D11")
$d11.NumberFormat = "@"
$d11.Value = '0.0723'
$d11.Style = "Normal"
$ws.Range("E11").Value = '  +4.23%  '

$d12 = $ws.Range("D12")
$d12.NumberFormat = "@"
$d12.Value = '0.0993'
$d12.Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").Value = '2.182.05'
$ws.Range("E13").Value = '  +2.42%  '

$d14 = $ws.Range("D14")
$d14.NumberFormat = "@"
$d14.Value = '12.05'
$d14.Style = "Normal"
$ws.Range("E14").Value = '  +4.96%  '

$d15 = $ws.Range("D15")
$d15.NumberFormat = "@"
$d15.Value = '0.698'
$d15.Style = "Normal"
$ws.Range("E15").Value = '  +3.09%  '

$ws.Range("D16").Value = '1.908.83'
$ws.Range("E16").Value = '  +2.48%  '

$d17 = $ws.Range("D17")
$d17.NumberFormat = "@"
$d17.Value = '4.85'
$d17.Style = "Normal"
$ws.Range("E17").Value = '  +2.77%  '

$ws.Range("D18").Value = '35.352.83'
$ws.Range("E18").Value = '  +0.51%  '

$d19 = $ws.Range("D19")
$d19.NumberFormat = "@"
$d19.Value = '72.18'
$d19.Style = "Normal"
$ws.Range("E19").Value = '  +3.23%  '

$ws.Range("D20").Value = '0.0₃0829'
$ws.Range("E20").Value = '  +4.34%  '

$d21 = $ws.Range("D21")
$d21.NumberFormat = "@"
$d21.Value = '240.42'
$d21.Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '

$d22 = $ws.Range("D22")
$d22.NumberFormat = "@"
$d22.Value = '12.47'
$d22.Style = "Normal"
$ws.Range("E22").Value = '  +1.98%  '

$ws.Range("E23").Value = '  +1.95%  '

$ws.Range("E24").Value = '  -0.49%  '

$ws.Range("E25").Value = '  +1.23%  '

$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = '2.36'
$d26.Style = "Normal"
$ws.Range("E26").Value = '  +25.51%  '

$d27 = $ws.Range("D27")
$d27.NumberFormat = "@"
$d27.Value = '170.14'
$d27.Style = "Normal"
$ws.Range("E27").Value = '  +0.46%  '

$d28 = $ws.Range("D28")
$d28.NumberFormat = "@"
$d28.Value = '8.42'
$d28.Style = "Normal"
$ws.Range("E28").Value = '  +4.51%  '

$d29 = $ws.Range("D29")
$d29.NumberFormat = "@"
$d29.Value = '18.39'
$d29.Style = "Normal"
$ws.Range("E29").Value = '  +3.77%  '

$ws.Range("E30").Value = '  +2.47%  '

$ws.Range("E31").Value = '  +3.14%  '

$d32 = $ws.Range("D32")
$d32.NumberFormat = "@"
$d32.Value = '0.0565'
$d32.Style = "Normal"
$ws.Range("E32").Value = '  +0.25%  '

$ws.Range("E33").Value = '  +0.51%  '

$d34 = $ws.Range("D34")
$d34.NumberFormat = "@"
$d34.Value = '0.929'
$d34.Style = "Normal"
$ws.Range("E34").Value = '  +13.32%  '

$ws.Range("E35").Value = '  +1.76%  '

$ws.Range("E36").Value = '  -4.57%  '

$ws.Range("E37").Value = '  -1.38%  '

$ws.Range("E38").Value = '  +2.26%  '

$d39 = $ws.Range("D39")
$d39.NumberFormat = "@"
$d39.Value = '1.10'
$d39.Style = "Normal"
$ws.Range("E39").Value = '  +0.36%  '

$ws.Range("E40").Value = '  +2.82%  '

$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = '16.27'
$d41.Style = "Normal"
$ws.Range("E41").Value = '  +7.79%  '

$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = '0.0628'
$d42.Style = "Normal"
$ws.Range("E42").Value = '  +4.91%  '

$d43 = $ws.Range("D43")
$d43.NumberFormat = "@"
$d43.Value = '89.86'
$d43.Style = "Normal"
$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("D44").Value = '1.338.39'
$ws.Range("E44").Value = '  -0.77%  '

$ws.Range("E45").Value = '  +2.96%  '

$d46 = $ws.Range("D46")
$d46.NumberFormat = "@"
$d46.Value = '48.10'
$d46.Style = "Normal"
$ws.Range("E46").Value = '  +38.88%  '

$ws.Range("E47").Value = '  -0.79%  '

$d48 = $ws.Range("D48")
$d48.NumberFormat = "@"
$d48.Value = '2.77'
$d48.Style = "Normal"
$ws.Range("E48").Value = '  +1.28%  '

$d49 = $ws.Range("D49")
$d49.NumberFormat = "@"
$d49.Value = '6.54'
$d49.Style = "Normal"
$ws.Range("E49").Value = '  -0.52%  '

$d50 = $ws.Range("D50")
$d50.NumberFormat = "@"
$d50.Value = '11.85'
$d50.Style = "Normal"
$ws.Range("E50").Value = '  -4.04%  '

$ws.Range("D51").Value = '2.092.63'
$ws.Range("E51").Value = '  +2.37%  '
